$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 31, duplicating the course data of ISYE 6740 (row 12) but
# flagged as the CS-required version of the course.
$ws.Range("A31").Value = "ISYE 6740 - CS Req"
$ws.Range("B31").Value = "CDA"
$ws.Range("C31").Value = "Computational Data Analytics"
$ws.Range("D31").Value = 3.63
$ws.Range("E31").Value = 14.6
$ws.Range("F31").Value = 4.03
$ws.Range("G31").Value = "CSE 6040 + Probability"
$ws.Range("I31").Value = "Python & MATLAB"
$ws.Range("J31").Value = "Jupyter Notebook"
$ws.Range("K31").Value = "Fall 2022, Summer 2022, Spring 2022, Fall 2021, Spring 2021"
$ws.Range("L31").Value = $true
$ws.Range("M31").Value = $true
$ws.Range("N31").Value = $true
$ws.Range("O31").Value = "CS Required"
